# Apply the "pushing changes to Timeline and Rh" edit.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.1.2"
$meta.Range("B5").Value = "CodeSystem - Transplant Timeline - NMDP"
$meta.Range("B8").Value = "2025-04-16T10:37:17-05:00"

# --- Concepts sheet: re-sort / re-code the timeline rows ---
# Column A ("Level") stays "1" for every data row in both the old and new
# layout, so it is left untouched (rewriting it would coerce the shared
# text string "1" into a numeric cell). Only Code (B) and Display (C) move.
$concepts = $wb.Worksheets.Item("Concepts")

$concepts.Cells.Item(2, 2).Value = "12W-6M"
$concepts.Cells.Item(2, 3).Value = "Over 12 weeks - up to 6 months"

$concepts.Cells.Item(3, 2).Value = "4-6W"
$concepts.Cells.Item(3, 3).Value = "Between 4-6 weeks"

$concepts.Cells.Item(4, 2).Value = "4W"
$concepts.Cells.Item(4, 3).Value = "Less than 4 weeks"

$concepts.Cells.Item(5, 2).Value = "6MG"
$concepts.Cells.Item(5, 3).Value = "Greater than 6 months"

$concepts.Cells.Item(6, 2).Value = "7-12w"
$concepts.Cells.Item(6, 3).Value = "Between 7-12 weeks"

$concepts.Cells.Item(7, 2).Value = "NA"
$concepts.Cells.Item(7, 3).Value = "N/A : MUD Transplant not preferred treatment"

$concepts.Cells.Item(8, 2).Value = "PEND"
$concepts.Cells.Item(8, 3).Value = "Pending, Case manager to follow up"
